$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.102.90"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "2.280.57"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.14"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.47"
$ws.Range("E6").Value = "  -3.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.601"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.81"
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0897"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.20"
$ws.Range("E12").Value = "  -2.73%  "
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").Value = "2.624.33"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "2.282.80"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "42.186.34"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.28"
$ws.Range("E19").Value = "  -3.95%  "
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.64"
$ws.Range("E21").Value = "  +28.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.52"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.47"
$ws.Range("E24").Value = "  +3.62%  "
$ws.Range("E25").Value = "  -4.95%  "
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.74"
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.32"
$ws.Range("E28").Value = "  +2.21%  "
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.17"
$ws.Range("E30").Value = "  +4.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "163.74"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.02"
$ws.Range("E32").Value = "  +2.74%  "
$ws.Range("E33").Value = "  -2.18%  "
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.114"
$ws.Range("E35").Value = "  -3.34%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.52"
$ws.Range("E36").Value = "  -13.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.53"
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0352"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.76"
$ws.Range("E39").Value = "  -4.96%  "
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "67.78"
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.86"
$ws.Range("E45").Value = "  -8.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "114.90"
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.85"
$ws.Range("E47").Value = "  -1.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "78.79"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("E49").Value = "  -2.15%  "
$ws.Range("D50").Value = "1.604.31"
$ws.Range("E50").Value = "  +4.02%  "
$ws.Range("E51").Value = "  -2.05%  "
